$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.042727768088382
$ws.Cells.Item(2, 4).Value = 1.048006095694486
$ws.Cells.Item(2, 5).Value = 1.046388445312651
$ws.Cells.Item(2, 6).Value = 1.057124699801312
$ws.Cells.Item(2, 9).Value = 1.040921202142737
$ws.Cells.Item(2, 10).Value = 1.047802032567034
$ws.Cells.Item(2, 11).Value = 1.05076714652594
$ws.Cells.Item(2, 12).Value = 1.049154020425594
$ws.Cells.Item(2, 13).Value = 1.059860528130246
$ws.Cells.Item(2, 14).Value = 1.01968522104141

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.044082252104491
$ws.Cells.Item(3, 4).Value = 1.049036399822995
$ws.Cells.Item(3, 5).Value = 1.047692616836325
$ws.Cells.Item(3, 6).Value = 1.058301770465155
$ws.Cells.Item(3, 9).Value = 1.04126433826696
$ws.Cells.Item(3, 10).Value = 1.048801062458678
$ws.Cells.Item(3, 11).Value = 1.051608843653737
$ws.Cells.Item(3, 12).Value = 1.050268544946509
$ws.Cells.Item(3, 13).Value = 1.060850449907363
$ws.Cells.Item(3, 14).Value = 1.020028347773158

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.044957888180578
$ws.Cells.Item(4, 4).Value = 1.049702145237498
$ws.Cells.Item(4, 5).Value = 1.048536061323979
$ws.Cells.Item(4, 6).Value = 1.059062745425296
$ws.Cells.Item(4, 9).Value = 1.041484329366623
$ws.Cells.Item(4, 10).Value = 1.049446249338465
$ws.Cells.Item(4, 11).Value = 1.052151945423013
$ws.Cells.Item(4, 12).Value = 1.050988731580584
$ws.Cells.Item(4, 13).Value = 1.061489752605553
$ws.Cells.Item(4, 14).Value = 1.020249664368951

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.045325817080077
$ws.Cells.Item(5, 4).Value = 1.049981804269934
$ws.Cells.Item(5, 5).Value = 1.048890543199007
$ws.Cells.Item(5, 6).Value = 1.059382502004516
$ws.Cells.Item(5, 9).Value = 1.041576326088189
$ws.Cells.Item(5, 10).Value = 1.049717189013697
$ws.Cells.Item(5, 11).Value = 1.052379900776649
$ws.Cells.Item(5, 12).Value = 1.051291265335152
$ws.Cells.Item(5, 13).Value = 1.061758220420542
$ws.Cells.Item(5, 14).Value = 1.020342536953464

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.045387583033795
$ws.Cells.Item(6, 4).Value = 1.050028747399625
$ws.Cells.Item(6, 5).Value = 1.048950056376339
$ws.Cells.Item(6, 6).Value = 1.05943618139666
$ws.Cells.Item(6, 9).Value = 1.041591744183876
$ws.Cells.Item(6, 10).Value = 1.04976266365182
$ws.Cells.Item(6, 11).Value = 1.052418154167093
$ws.Cells.Item(6, 12).Value = 1.051342048529489
$ws.Cells.Item(6, 13).Value = 1.061803280119407
$ws.Cells.Item(6, 14).Value = 1.020358120790212

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.044962805199819
$ws.Cells.Item(7, 4).Value = 1.049705882918044
$ws.Cells.Item(7, 5).Value = 1.048540798325257
$ws.Cells.Item(7, 6).Value = 1.059067018643028
$ws.Cells.Item(7, 9).Value = 1.041485560545703
$ws.Cells.Item(7, 10).Value = 1.049449870810641
$ws.Cells.Item(7, 11).Value = 1.052154992803558
$ws.Cells.Item(7, 12).Value = 1.050992774961874
$ws.Cells.Item(7, 13).Value = 1.061493341041841
$ws.Cells.Item(7, 14).Value = 1.020250905999556

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.043185691753339
$ws.Cells.Item(8, 4).Value = 1.048354485841068
$ws.Cells.Item(8, 5).Value = 1.046829289785686
$ws.Cells.Item(8, 6).Value = 1.057522636407154
$ws.Cells.Item(8, 9).Value = 1.041037589830685
$ws.Cells.Item(8, 10).Value = 1.048139920438766
$ws.Cells.Item(8, 11).Value = 1.051051920705431
$ws.Cells.Item(8, 12).Value = 1.04953088511974
$ws.Cells.Item(8, 13).Value = 1.060195335692854
$ws.Cells.Item(8, 14).Value = 1.019801329747975

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.040047824282761
$ws.Cells.Item(9, 4).Value = 1.045965920149088
$ws.Cells.Item(9, 5).Value = 1.043809850322766
$ws.Cells.Item(9, 6).Value = 1.054795994402777
$ws.Cells.Item(9, 9).Value = 1.040232529014877
$ws.Cells.Item(9, 10).Value = 1.045821898895893
$ws.Cells.Item(9, 11).Value = 1.049096326903524
$ws.Cells.Item(9, 12).Value = 1.046947160084611
$ws.Cells.Item(9, 13).Value = 1.057898452162483
$ws.Cells.Item(9, 14).Value = 1.019003645390949

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.037951355787452
$ws.Cells.Item(10, 4).Value = 1.044368528942411
$ws.Cells.Item(10, 5).Value = 1.041794282774575
$ws.Cells.Item(10, 6).Value = 1.052974532246122
$ws.Cells.Item(10, 9).Value = 1.039685210911379
$ws.Cells.Item(10, 10).Value = 1.044269829825024
$ws.Cells.Item(10, 11).Value = 1.047784491808407
$ws.Cells.Item(10, 12).Value = 1.045219306001759
$ws.Cells.Item(10, 13).Value = 1.056360559281452
$ws.Cells.Item(10, 14).Value = 1.018468117664538

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.037042418834361
$ws.Cells.Item(11, 4).Value = 1.043675616928875
$ws.Cells.Item(11, 5).Value = 1.040920850193929
$ws.Cells.Item(11, 6).Value = 1.052184901278886
$ws.Cells.Item(11, 9).Value = 1.039445683450738
$ws.Cells.Item(11, 10).Value = 1.043596130677765
$ws.Cells.Item(11, 11).Value = 1.047214496997537
$ws.Cells.Item(11, 12).Value = 1.044469806690052
$ws.Cells.Item(11, 13).Value = 1.055693023785115
$ws.Cells.Item(11, 14).Value = 1.018235328393339

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.036704620624987
$ws.Cells.Item(12, 4).Value = 1.043418050503407
$ws.Cells.Item(12, 5).Value = 1.040596312013957
$ws.Cells.Item(12, 6).Value = 1.051891454534432
$ws.Cells.Item(12, 9).Value = 1.03935632996433
$ws.Cells.Item(12, 10).Value = 1.043345638273905
$ws.Cells.Item(12, 11).Value = 1.04700247789046
$ws.Cells.Item(12, 12).Value = 1.044191205874916
$ws.Cells.Item(12, 13).Value = 1.055444825072325
$ws.Cells.Item(12, 14).Value = 1.018148723362786

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.036777087655226
$ws.Cells.Item(13, 4).Value = 1.043473307935701
$ws.Cells.Item(13, 5).Value = 1.040665931443765
$ws.Cells.Item(13, 6).Value = 1.051954406369484
$ws.Cells.Item(13, 9).Value = 1.039375513914811
$ws.Cells.Item(13, 10).Value = 1.043399381142192
$ws.Cells.Item(13, 11).Value = 1.047047970225823
$ws.Cells.Item(13, 12).Value = 1.044250975964534
$ws.Cells.Item(13, 13).Value = 1.055498075737302
$ws.Cells.Item(13, 14).Value = 1.018167306642191

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.037014500003032
$ws.Cells.Item(14, 4).Value = 1.043654330261989
$ws.Cells.Item(14, 5).Value = 1.040894025963927
$ws.Cells.Item(14, 6).Value = 1.05216064780998
$ws.Cells.Item(14, 9).Value = 1.039438305271379
$ws.Cells.Item(14, 10).Value = 1.043575430030043
$ws.Cells.Item(14, 11).Value = 1.047196977521245
$ws.Cells.Item(14, 12).Value = 1.044446781632381
$ws.Cells.Item(14, 13).Value = 1.055672512645638
$ws.Cells.Item(14, 14).Value = 1.018228172390606

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.037160753762639
$ws.Cells.Item(15, 4).Value = 1.043765839080005
$ws.Cells.Item(15, 5).Value = 1.041034548288146
$ws.Cells.Item(15, 6).Value = 1.052287700959862
$ws.Cells.Item(15, 9).Value = 1.039476942393471
$ws.Cells.Item(15, 10).Value = 1.043683866249517
$ws.Cells.Item(15, 11).Value = 1.04728874630288
$ws.Cells.Item(15, 12).Value = 1.044567396893718
$ws.Cells.Item(15, 13).Value = 1.055779956245905
$ws.Cells.Item(15, 14).Value = 1.018265655633506

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.038011654199025
$ws.Cells.Item(16, 4).Value = 1.04441448902956
$ws.Cells.Item(16, 5).Value = 1.041852234940016
$ws.Cells.Item(16, 6).Value = 1.053026917644255
$ws.Cells.Item(16, 9).Value = 1.039701054024661
$ws.Cells.Item(16, 10).Value = 1.044314506087106
$ws.Cells.Item(16, 11).Value = 1.047822278907409
$ws.Cells.Item(16, 12).Value = 1.045269019499906
$ws.Cells.Item(16, 13).Value = 1.056404827052286
$ws.Cells.Item(16, 14).Value = 1.018483548012248

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.038545089237457
$ws.Cells.Item(17, 4).Value = 1.044821038232197
$ws.Cells.Item(17, 5).Value = 1.042364963350933
$ws.Cells.Item(17, 6).Value = 1.053490358370702
$ws.Cells.Item(17, 9).Value = 1.039840953588346
$ws.Cells.Item(17, 10).Value = 1.044709647311991
$ws.Cells.Item(17, 11).Value = 1.048156422571097
$ws.Cells.Item(17, 12).Value = 1.045708770694145
$ws.Cells.Item(17, 13).Value = 1.056796356592193
$ws.Cells.Item(17, 14).Value = 1.018619983760751

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.038856122063905
$ws.Cells.Item(18, 4).Value = 1.045058052933854
$ws.Cells.Item(18, 5).Value = 1.042663963995444
$ws.Cells.Item(18, 6).Value = 1.053760586418257
$ws.Cells.Item(18, 9).Value = 1.039922310134761
$ws.Cells.Item(18, 10).Value = 1.044939968229922
$ws.Cells.Item(18, 11).Value = 1.048351133856447
$ws.Cells.Item(18, 12).Value = 1.045965142354288
$ws.Cells.Item(18, 13).Value = 1.057024573283664
$ws.Cells.Item(18, 14).Value = 1.018699477452728

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.038962157588975
$ws.Cells.Item(19, 4).Value = 1.045138848803478
$ws.Cells.Item(19, 5).Value = 1.042765904570558
$ws.Cells.Item(19, 6).Value = 1.05385271219493
$ws.Cells.Item(19, 9).Value = 1.039950009186332
$ws.Cells.Item(19, 10).Value = 1.045018474976315
$ws.Cells.Item(19, 11).Value = 1.048417493393711
$ws.Cells.Item(19, 12).Value = 1.046052536913731
$ws.Cells.Item(19, 13).Value = 1.057102362938331
$ws.Cells.Item(19, 14).Value = 1.018726568030909

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.03848786820835
$ws.Cells.Item(20, 4).Value = 1.04477743163807
$ws.Cells.Item(20, 5).Value = 1.042309959218832
$ws.Cells.Item(20, 6).Value = 1.053440644808367
$ws.Cells.Item(20, 9).Value = 1.039825968990221
$ws.Cells.Item(20, 10).Value = 1.044667268797053
$ws.Cells.Item(20, 11).Value = 1.048120591696083
$ws.Cells.Item(20, 12).Value = 1.045661602799112
$ws.Cells.Item(20, 13).Value = 1.056754365324542
$ws.Cells.Item(20, 14).Value = 1.018605354501012

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.036944592956773
$ws.Cells.Item(21, 4).Value = 1.043601028914892
$ws.Cells.Item(21, 5).Value = 1.040826860787015
$ws.Cells.Item(21, 6).Value = 1.052099918812514
$ws.Cells.Item(21, 9).Value = 1.039419825347592
$ws.Cells.Item(21, 10).Value = 1.043523594966568
$ws.Cells.Item(21, 11).Value = 1.047153106827223
$ws.Cells.Item(21, 12).Value = 1.044389127390744
$ws.Cells.Item(21, 13).Value = 1.055621152149464
$ws.Cells.Item(21, 14).Value = 1.018210252728678

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.035973238776411
$ws.Cells.Item(22, 4).Value = 1.042860288271819
$ws.Cells.Item(22, 5).Value = 1.03989375992609
$ws.Cells.Item(22, 6).Value = 1.051256124090332
$ws.Cells.Item(22, 9).Value = 1.039162253991785
$ws.Cells.Item(22, 10).Value = 1.042803070415905
$ws.Cells.Item(22, 11).Value = 1.046543087070916
$ws.Cells.Item(22, 12).Value = 1.043587892783454
$ws.Cells.Item(22, 13).Value = 1.054907229358055
$ws.Cells.Item(22, 14).Value = 1.017961044893762

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.036488272009132
$ws.Cells.Item(23, 4).Value = 1.043253073218876
$ws.Cells.Item(23, 5).Value = 1.040388474432551
$ws.Cells.Item(23, 6).Value = 1.051703515225689
$ws.Cells.Item(23, 9).Value = 1.039299007669271
$ws.Cells.Item(23, 10).Value = 1.043185172932456
$ws.Cells.Item(23, 11).Value = 1.046866634456044
$ws.Cells.Item(23, 12).Value = 1.044012755494061
$ws.Cells.Item(23, 13).Value = 1.055285829665999
$ws.Cells.Item(23, 14).Value = 1.018093230119039

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.03851372427041
$ws.Cells.Item(24, 4).Value = 1.044797135947185
$ws.Cells.Item(24, 5).Value = 1.042334813419628
$ws.Cells.Item(24, 6).Value = 1.053463108502424
$ws.Cells.Item(24, 9).Value = 1.039832740640522
$ws.Cells.Item(24, 10).Value = 1.044686418312644
$ws.Cells.Item(24, 11).Value = 1.048136782710215
$ws.Cells.Item(24, 12).Value = 1.045682916335151
$ws.Cells.Item(24, 13).Value = 1.056773339852362
$ws.Cells.Item(24, 14).Value = 1.018611965102463

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.040859821069096
$ws.Cells.Item(25, 4).Value = 1.046584294172724
$ws.Cells.Item(25, 5).Value = 1.044590890815755
$ws.Cells.Item(25, 6).Value = 1.055501536060084
$ws.Cells.Item(25, 9).Value = 1.04044252166813
$ws.Cells.Item(25, 10).Value = 1.046422334074884
$ws.Cells.Item(25, 11).Value = 1.049603312763479
$ws.Cells.Item(25, 12).Value = 1.047616046717809
$ws.Cells.Item(25, 13).Value = 1.058493409584187
$ws.Cells.Item(25, 14).Value = 1.019210520515105

